$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.414.46"
$ws.Range("E2").Value = "  -3.30%  "

$ws.Range("D3").Value = "3.130.09"
$ws.Range("E3").Value = "  -4.21%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "561.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.43%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.582"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -8.52%  "

$ws.Range("D9").Value = "3.118.07"
$ws.Range("E9").Value = "  -4.55%  "

$ws.Range("E10").Value = "  -1.99%  "

$ws.Range("E11").Value = "  -6.94%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.378"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.58%  "

$ws.Range("D13").Value = "3.665.06"
$ws.Range("E13").Value = "  -4.57%  "

$ws.Range("E14").Value = "  -1.13%  "

$ws.Range("D15").Value = "63.420.11"
$ws.Range("E15").Value = "  -3.46%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "24.80"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.96%  "

$ws.Range("D17").Value = "3.125.17"
$ws.Range("E17").Value = "  -4.91%  "

$ws.Range("E18").Value = "  -5.49%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "400.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.62%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.38%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.65%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.90%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.05%  "

$ws.Range("E25").Value = "  -2.65%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.477"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.98%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000100"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -10.47%  "

$ws.Range("E28").Value = "  -7.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.11%  "

$ws.Range("E30").Value = "  -0.04%  "

$ws.Range("E31").Value = "  -6.53%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.86"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.47%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.57%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.75"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.49%  "

$ws.Range("E35").Value = "  -6.03%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "152.57"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.68%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.32"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.86%  "

$ws.Range("D38").Value = "2.733.98"
$ws.Range("E38").Value = "  -3.64%  "

$ws.Range("E39").Value = "  -7.32%  "

$ws.Range("E40").Value = "  -5.85%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.19"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -10.25%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.54"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.38%  "

$ws.Range("E43").Value = "  -6.88%  "

$ws.Range("E44").Value = "  -2.97%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.34"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.10%  "

$ws.Range("E46").Value = "  -3.76%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "20.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -8.14%  "

$ws.Range("E48").Value = "  -0.07%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "279.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -9.84%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0972"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.12%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.91%  "
